# Actualización automática 2025-06-01 08:00:06
#
# This script reproduces, via Excel COM-interop, the monthly "roll forward"
# update applied to the workbook:
#   - Sheet "VENTA MENSUAL": the oldest month column is dropped, every
#     month's figures shift one column to the left (C<-D, D<-E, E<-F),
#     a new (empty) month is appended in column F, and the month name
#     headers + a few column widths shift accordingly.
#   - Sheet "VENTAS POR GRUPO": the figures that belonged exclusively to
#     the now-dropped month are zeroed out (both the numeric cells and
#     the "n de 54" counter labels on the totals row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" -- zero out cells whose only contribution
# was from the month that just rolled off the report.
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$zeroCells = @("L3","N3","C10","L10","D12","L12","L16","C26","K26","L26","L27","N27","J29","K29","L42","L43","D46","E46","L47","N47","L55")
foreach ($ref in $zeroCells) {
    $wsGrupo.Range($ref).Value = 0
}

$zeroCounters = @("C56","D56","E56","J56","K56","L56","N56")
foreach ($ref in $zeroCounters) {
    $wsGrupo.Range($ref).Value = "0 de 54"
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" -- shift every month's data one column to the
# left and append a fresh, empty month on the right.
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Month name headers (row 1, columns C:F) shift left; a new month name
# ("junio") is introduced in the rightmost column.
$wsMensual.Cells.Item(1, 3).Value = $wsMensual.Cells.Item(1, 4).Value()
$wsMensual.Cells.Item(1, 4).Value = $wsMensual.Cells.Item(1, 5).Value()
$wsMensual.Cells.Item(1, 5).Value = $wsMensual.Cells.Item(1, 6).Value()
$wsMensual.Cells.Item(1, 6).Value = "junio"

# Data rows (2:56, including the totals row) shift left the same way,
# with the new rightmost month starting at 0.
for ($r = 2; $r -le 56; $r++) {
    $dVal = $wsMensual.Cells.Item($r, 4).Value()
    $eVal = $wsMensual.Cells.Item($r, 5).Value()
    $fVal = $wsMensual.Cells.Item($r, 6).Value()

    $wsMensual.Cells.Item($r, 3).Value = $dVal
    $wsMensual.Cells.Item($r, 4).Value = $eVal
    $wsMensual.Cells.Item($r, 5).Value = $fVal
    $wsMensual.Cells.Item($r, 6).Value = 0
}

# Column widths follow the same leftward shift (min/max="3" column takes
# on the former width of column 4, etc.); column F gets a brand-new,
# narrower default width. ColumnWidth is expressed in character units,
# which Excel stores in the XML offset by 5/6 of a character from the
# value you assign -- correct for that so the saved <col width="..."/>
# matches exactly.
$pad = 0.8333333333333334
$wsMensual.Columns.Item(3).ColumnWidth = 13 - $pad
$wsMensual.Columns.Item(5).ColumnWidth = 14 - $pad
$wsMensual.Columns.Item(6).ColumnWidth = 11 - $pad
